$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for columns E-J (Clear, Assertive, Cautious, Optimistic, Specific, Relevant)
# for rows 2..32
$data = @{
    2  = @(2,2,1,2,2,2)
    3  = @(2,2,1,1,1,2)
    4  = @(2,2,1,1,2,2)
    5  = @(2,2,1,2,2,2)
    6  = @(2,2,1,2,2,2)
    7  = @(2,1,1,1,1,2)
    8  = @(1,2,1,2,1,1)
    9  = @(2,2,1,2,2,2)
    10 = @(2,2,1,2,2,2)
    11 = @(2,2,1,2,2,2)
    12 = @(2,2,1,2,1,2)
    13 = @(2,2,1,2,2,2)
    14 = @(2,2,1,2,2,2)
    15 = @(1,1,1,1,1,1)
    16 = @(2,2,1,2,2,2)
    17 = @(2,2,2,2,2,2)
    18 = @(2,2,1,2,2,2)
    19 = @(2,2,2,2,2,2)
    20 = @(1,1,1,2,1,2)
    21 = @(2,1,0,1,1,2)
    22 = @(2,2,1,1,2,2)
    23 = @(2,2,2,2,2,2)
    24 = @(2,2,1,2,2,2)
    25 = @(2,2,1,2,2,2)
    26 = @(2,2,1,1,2,2)
    27 = @(2,1,2,2,1,2)
    28 = @(2,1,1,1,1,2)
    29 = @(2,2,1,1,1,2)
    30 = @(2,2,2,2,2,2)
    31 = @(2,2,1,1,2,2)
    32 = @(2,2,1,2,1,2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 5 + $i  # E=5
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}

# Update sheet view: zoom, freeze panes (header row frozen), and final selection
$ws.Select()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.Zoom = 85
$ws.Range("F33").Select()
